$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1.01
$ws.Range("F3").Value = 1.78
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 4.9
$ws.Range("I3").Value = 5.2
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 2.22
$ws.Range("V3").Value = 1.24
$ws.Range("W3").Value = 2.24
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 980
$ws.Range("AA3").Value = 120
$ws.Range("AD3").Value = 19
$ws.Range("AE3").Value = 60
$ws.Range("AI3").Value = 60
$ws.Range("AJ3").Value = 19
$ws.Range("AM3").Value = 90
$ws.Range("AN3").Value = 9.4
$ws.Range("G4").Value = 5.8
$ws.Range("L4").Value = 1.54
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.62
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 1.55
$ws.Range("Q4").Value = 2.46
$ws.Range("R4").Value = 1.19
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = 2.12
$ws.Range("U4").Value = 1.71
$ws.Range("V4").Value = 1.93
$ws.Range("W4").Value = 1.21
$ws.Range("AB4").Value = 14
$ws.Range("AD4").Value = 11.5
$ws.Range("AE4").Value = 28
$ws.Range("AG4").Value = 26
$ws.Range("AJ4").Value = 190
$ws.Range("AK4").Value = 120
$ws.Range("AL4").Value = 140
$ws.Range("AM4").Value = 250
$ws.Range("AN4").Value = 210
$ws.Range("AO4").Value = 25
$ws.Range("H5").Value = 1.56
$ws.Range("I5").Value = 1.6
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.75
$ws.Range("P5").Value = 1.93
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.35
$ws.Range("S5").Value = 3.5
$ws.Range("T5").Value = 2.06
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 2.62
$ws.Range("W5").Value = 1.15
$ws.Range("X5").Value = 18.5
$ws.Range("AG5").Value = 980
$ws.Range("F6").Value = 4.3
$ws.Range("G6").Value = 5.3
$ws.Range("H6").Value = 1.92
$ws.Range("I6").Value = 2.04
$ws.Range("K6").Value = 3.9
$ws.Range("V6").Value = 1.96
$ws.Range("W6").Value = 1.24
$ws.Range("AD6").Value = 1000
$ws.Range("G7").Value = 1.29
$ws.Range("H7").Value = 11.5
$ws.Range("J7").Value = 6.6
$ws.Range("S7").Value = 2.14
$ws.Range("W7").Value = 4.4
$ws.Range("P8").Value = 2.1
$ws.Range("U8").Value = 2.26
$ws.Range("V8").Value = 1.31
$ws.Range("Y8").Value = 1000
$ws.Range("Z8").Value = 34
$ws.Range("AF8").Value = 1000
$ws.Range("AH8").Value = 980
$ws.Range("F9").Value = 3.6
$ws.Range("G9").Value = 3.75
$ws.Range("H9").Value = 2.26
$ws.Range("I9").Value = 2.28
$ws.Range("K9").Value = 3.55
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 2.04
$ws.Range("Q9").Value = 1.89
$ws.Range("V9").Value = 1.78
$ws.Range("W9").Value = 1.37
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 15
$ws.Range("AA9").Value = 29
$ws.Range("AG9").Value = 15
$ws.Range("AJ9").Value = 70
$ws.Range("AK9").Value = 42
$ws.Range("AL9").Value = 55
$ws.Range("AM9").Value = 85
$ws.Range("AN9").Value = 38
$ws.Range("AO9").Value = 19.5
$ws.Range("H10").Value = 1.85
$ws.Range("I10").Value = 2.02
$ws.Range("J10").Value = 3.45
$ws.Range("K10").Value = 4.2
$ws.Range("Q10").Value = 1.63
$ws.Range("V10").Value = 1.98
$ws.Range("F11").Value = 1.23
$ws.Range("G11").Value = 1.27
$ws.Range("H11").Value = 14
$ws.Range("I11").Value = 17.5
$ws.Range("J11").Value = 6.6
$ws.Range("K11").Value = 7.8
$ws.Range("L11").Value = 1.22
$ws.Range("N11").Value = 5.4
$ws.Range("O11").Value = 1.18
$ws.Range("P11").Value = 2.52
$ws.Range("Q11").Value = 1.54
$ws.Range("R11").Value = 1.61
$ws.Range("T11").Value = 2.14
$ws.Range("U11").Value = 1.72
$ws.Range("W11").Value = 4.4
$ws.Range("Z11").Value = 190
$ws.Range("AB11").Value = 10.5
$ws.Range("AC11").Value = 17
$ws.Range("AD11").Value = 60
$ws.Range("AE11").Value = 330
$ws.Range("AF11").Value = 8.6
$ws.Range("AG11").Value = 14.5
$ws.Range("AI11").Value = 230
$ws.Range("AJ11").Value = 11.5
$ws.Range("AL11").Value = 44
$ws.Range("AM11").Value = 240
$ws.Range("AN11").Value = 4.3
$ws.Range("F12").Value = 2.52
$ws.Range("G12").Value = 2.7
$ws.Range("Q12").Value = 1.92
$ws.Range("V12").Value = 1.47
$ws.Range("W12").Value = 1.59
$ws.Range("AA12").Value = 60
$ws.Range("AC12").Value = 9.6
$ws.Range("AJ12").Value = 40
$ws.Range("AL12").Value = 46
$ws.Range("F13").Value = 1.47
$ws.Range("G13").Value = 1.58
$ws.Range("I13").Value = 12.5
$ws.Range("K13").Value = 4.6
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 2.86
$ws.Range("O13").Value = 1.43
$ws.Range("P13").Value = 1.63
$ws.Range("Q13").Value = 2.26
$ws.Range("R13").Value = 1.23
$ws.Range("S13").Value = 3.95
$ws.Range("T13").Value = 2.4
$ws.Range("U13").Value = 1.58
$ws.Range("V13").Value = 1.09
$ws.Range("W13").Value = 2.74
$ws.Range("X13").Value = 13.5
$ws.Range("AB13").Value = 7
$ws.Range("AC13").Value = 12
$ws.Range("AF13").Value = 9
$ws.Range("AJ13").Value = 16.5
$ws.Range("AN13").Value = 14.5
$ws.Range("F14").Value = 2.32
$ws.Range("G14").Value = 2.54
$ws.Range("H14").Value = 2.86
$ws.Range("K14").Value = 4.2
$ws.Range("V14").Value = 1.48
$ws.Range("W14").Value = 1.65
$ws.Range("AE14").Value = 980
$ws.Range("AM14").Value = 65
$ws.Range("F15").Value = 11.5
$ws.Range("J15").Value = 5.2
$ws.Range("L15").Value = 1.33
$ws.Range("N15").Value = 3.8
$ws.Range("P15").Value = 1.92
$ws.Range("Q15").Value = 1.99
$ws.Range("S15").Value = 3.5
$ws.Range("T15").Value = 2.42
$ws.Range("U15").Value = 1.64
$ws.Range("V15").Value = 3.55
$ws.Range("Y15").Value = 7.4
$ws.Range("AA15").Value = 10.5
$ws.Range("AE15").Value = 16.5
$ws.Range("AG15").Value = 48
$ws.Range("AH15").Value = 40
$ws.Range("AI15").Value = 55
$ws.Range("AJ15").Value = 710
$ws.Range("AL15").Value = 240
$ws.Range("AM15").Value = 290
$ws.Range("AN15").Value = 510
$ws.Range("AO15").Value = 7.2
$ws.Range("F16").Value = 4.3
$ws.Range("N16").Value = 6.2
$ws.Range("P16").Value = 2.86
$ws.Range("R16").Value = 1.72
$ws.Range("U16").Value = 2.58
$ws.Range("W16").Value = 1.27
$ws.Range("Y16").Value = 16
$ws.Range("Z16").Value = 17.5
$ws.Range("AG16").Value = 19.5
$ws.Range("AI16").Value = 26
$ws.Range("G17").Value = 2.78
$ws.Range("I17").Value = 2.66
$ws.Range("T17").Value = 1.55
$ws.Range("V17").Value = 1.6
$ws.Range("W17").Value = 1.56
$ws.Range("AF17").Value = 25
$ws.Range("Q18").Value = 2.22
$ws.Range("U18").Value = 1.75
$ws.Range("V18").Value = 1.14
$ws.Range("F19").Value = 2.56
$ws.Range("G19").Value = 2.58
$ws.Range("H19").Value = 3.15
$ws.Range("I19").Value = 3.2
$ws.Range("N19").Value = 3.6
$ws.Range("P19").Value = 1.86
$ws.Range("T19").Value = 1.84
$ws.Range("U19").Value = 2.14
$ws.Range("W19").Value = 1.63
$ws.Range("X19").Value = 12.5
$ws.Range("AA19").Value = 50
$ws.Range("AD19").Value = 14
$ws.Range("AF19").Value = 15
$ws.Range("AM19").Value = 85
$ws.Range("F20").Value = 1.84
$ws.Range("G20").Value = 1.85
$ws.Range("K20").Value = 3.65
$ws.Range("L20").Value = 1.54
$ws.Range("P20").Value = 1.67
$ws.Range("V20").Value = 1.21
$ws.Range("W20").Value = 2.16
$ws.Range("Y20").Value = 15
$ws.Range("AB20").Value = 6.8
$ws.Range("AL20").Value = 55
$ws.Range("P21").Value = 2.24
$ws.Range("Q21").Value = 1.74
$ws.Range("U21").Value = 1.76
$ws.Range("X21").Value = 22
$ws.Range("AA21").Value = 550
$ws.Range("AD21").Value = 42
$ws.Range("AH21").Value = 32
$ws.Range("AL21").Value = 40
$ws.Range("AM21").Value = 200
$ws.Range("AO21").Value = 300
